# Update the "Initial Time" policy-schedule input: advance the simulation's
# first year (IT sheet, cell B2) from 2017 to 2019.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IT")
$ws.Range("B2").Value = 2019
